$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values (shared strings will be rewritten to: username, quandohong28, zthanh13, admin)
$ws.Range("A2").Value = "quandohong28"
$ws.Range("A3").Value = "zthanh13"
$ws.Range("A4").Value = "admin"

# Update the active selection from B7 to D6
$ws.Range("D6").Select()
